$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.111.83'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '3.573.12'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.38'
$ws.Range("E5").Value = '  -1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '657.84'
$ws.Range("E6").Value = '  +2.99%  '
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  -1.61%  '
$ws.Range("D11").Value = '3.571.67'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.37'
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.45'
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").Value = '4.235.37'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("D16").Value = '95.020.93'
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").Value = '3.569.93'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("E19").Value = '  -3.39%  '
$ws.Range("E20").Value = '  -5.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.82'
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.45'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.47'
$ws.Range("E23").Value = '  -1.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.478'
$ws.Range("E24").Value = '  -3.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.85'
$ws.Range("E25").Value = '  +2.16%  '
$ws.Range("E26").Value = '  -1.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.15'
$ws.Range("E27").Value = '  -2.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.66'
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("D29").Value = '3.764.25'
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.04'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.143'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.52'
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  -2.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.90'
$ws.Range("E36").Value = '  +4.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.66'
$ws.Range("E37").Value = '  +11.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.557'
$ws.Range("E38").Value = '  -2.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.43'
$ws.Range("E39").Value = '  +6.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '579.75'
$ws.Range("E40").Value = '  +0.97%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.907'
$ws.Range("E43").Value = '  -1.86%  '
$ws.Range("E44").Value = '  +4.40%  '
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '34.52'
$ws.Range("E47").Value = '  +3.43%  '
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0414'
$ws.Range("E49").Value = '  -3.84%  '
$ws.Range("E50").Value = '  +1.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.17'
$ws.Range("E51").Value = '  +0.54%  '
